$wb = $excel.ActiveWorkbook

# The F column ("想去人数" / want-to-go count) was refreshed with newer
# scraped numbers for both the "展览" sheet and the aggregated "全部类型"
# sheet (they mirror the same rows).
$sheetNames = @("展览", "全部类型")

$updates = @{
    "F7"  = 47
    "F11" = 4677
    "F12" = 4476
    "F16" = 159
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
